$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 0.6648627278815172
$ws.Range("C4").Value = 0.667
$ws.Range("D4").Value = 0.677166889836075
$ws.Range("E4").Value = 0.6719999999999999
$ws.Range("F4").Value = 0.5240434850327348
$ws.Range("G4").Value = 0.53
$ws.Range("H4").Value = 0.5204935166877871
$ws.Range("I4").Value = 0.5205
$ws.Range("J4").Value = 0.6434232903758921
$ws.Range("K4").Value = 0.675
$ws.Range("L4").Value = 0.6209477564719348
$ws.Range("M4").Value = 0.6295000000000001

$ws.Range("B5").Value = 0.4008757564102146
$ws.Range("C5").Value = 0.351
$ws.Range("D5").Value = 0.5758951126392987
$ws.Range("E5").Value = 0.5945
$ws.Range("F5").Value = 0.6680529450303134
$ws.Range("H5").Value = 0.5099724786095357
$ws.Range("I5").Value = 0.5175000000000001
$ws.Range("J5").Value = 0.4035272835243034
$ws.Range("K5").Value = 0.4029999999999999
$ws.Range("L5").Value = 0.5745105902810784
$ws.Range("M5").Value = 0.5669999999999999

$ws.Range("B6").Value = 0.7218838531824505
$ws.Range("C6").Value = 0.713
$ws.Range("D6").Value = 0.7571914438514239
$ws.Range("E6").Value = 0.7375
$ws.Range("F6").Value = 0.5054654098681046
$ws.Range("G6").Value = 0.506
$ws.Range("H6").Value = 0.5074582709594921
$ws.Range("I6").Value = 0.5054999999999999
$ws.Range("J6").Value = 0.6999827819334058
$ws.Range("K6").Value = 0.6910000000000001
$ws.Range("L6").Value = 0.7396124781829192
$ws.Range("M6").Value = 0.7209999999999999
